$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 20240621
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 7

# Row 8
$ws.Range("B8").Value = "6,31"

# Row 9
$ws.Range("A9").Value = 20240718
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 8

# Row 10
$ws.Range("D10").Value = "9,10,11,16,17,18"

# Row 11
$ws.Range("A11").Value = 20240719
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 6

# Row 12
$ws.Range("E12").Value = 8

# Row 13
$ws.Range("A13").Value = 20240720
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 5

# Row 14
$ws.Range("A14").Value = "due to naming convention, I split the two sessions done on 20240719 into two separate days, otherwise that session date would contain double the amount of data which seems bad (need to maintain consistent trials)"

# Row 15
$ws.Range("A15").Value = 20240808
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 6

# Row 17 (row 16 intentionally left blank)
$ws.Range("A17").Value = 20240809
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 4

# Column widths
# NOTE: the host engine stores ColumnWidth as (value + 5/6) quantized to the
# nearest 1/6, so the values below are chosen to land as close as possible
# on the target stored widths (17.42578125, 20, 22.85546875, 21.85546875,
# 21.7109375).
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 21
$ws.Columns.Item(6).ColumnWidth = 20.833333333333332

# Selection
$ws.Range("C17").Select()
